$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.889471650123596
$ws.Range("B1").Value = -1
$ws.Range("C1").Value = 2.693895101547241
$ws.Range("D1").Value = 1.252211213111877
$ws.Range("E1").Value = 0.912026584148407
